# Daily attendance processing - 2025-10-12 06:25:19
# Normalizes the "Recorded By" (column G) values: the list of recorders is
# right-rotated by one element so that the entry which used to start the
# list ("System"/"system") moves down, with the former last entry now
# appearing first.
#   "System, dnasr281@gmail.com"                 -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com"        -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = @($val -split ",\s*")
        if ($parts.Count -gt 1 -and $parts[0].ToLower() -eq "system") {
            $rotated = @($parts[-1]) + @($parts[0..($parts.Count - 2)])
            $cell.Value2 = [string]::Join(", ", $rotated)
        }
    }
}
